# [LAB6] - Added drop down menu for machine type - Updated presentation
#
# On the "Demo" slide (slide 5), italicize the three paragraphs that show
# the Visual Studio Code / dotnet CLI instructions:
#   - "Visual Studio Code"
#   - "dotnet new mvc"
#   - "dotnet new xunit"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraphs (1-based): 5 = "Visual Studio Code", 6 = "dotnet new mvc",
# 7 = "dotnet new xunit". Set each paragraph's runs to italic.
$tr.Paragraphs(5, 1).Font.Italic = -1
$tr.Paragraphs(6, 1).Font.Italic = -1
$tr.Paragraphs(7, 1).Font.Italic = -1
